# VAN-1811: Prepare and write FUNCTIONAL test cases and test scripts
#
# Refresh the "Order Assign" test-data row with a new set of FuncLoc /
# SAID / Previous-Doc identifiers, and stage the newly prepared pool of
# test values (10 rows x 3 columns: Previous Doc / FuncLoc / SAID) on the
# spare Sheet2 for the upcoming functional test scripts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Stage the newly prepared test-data values -----------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A1").Value = "'7099614533"
$ws2.Range("B1").Value = "ABCD142660"
$ws2.Range("C1").Value = "'1336650338"

$ws2.Range("A2").Value = "'1785657360"
$ws2.Range("B2").Value = "ABCD400484"
$ws2.Range("C2").Value = "'4323691594"

$ws2.Range("A3").Value = "'7053075709"
$ws2.Range("B3").Value = "ABCD508618"
$ws2.Range("C3").Value = "'0610774080"

$ws2.Range("A4").Value = "'5326525470"
$ws2.Range("B4").Value = "ABCD647150"
$ws2.Range("C4").Value = "'2010233511"

$ws2.Range("A5").Value = "'2838288423"
$ws2.Range("B5").Value = "ABCD046558"
$ws2.Range("C5").Value = "'1366406344"

$ws2.Range("A6").Value = "'1859487456"
$ws2.Range("B6").Value = "ABCD833949"
$ws2.Range("C6").Value = "'0387558984"

$ws2.Range("A7").Value = "'1769718942"
$ws2.Range("B7").Value = "ABCD614732"
$ws2.Range("C7").Value = "'9472362756"

$ws2.Range("A8").Value = "'5370729632"
$ws2.Range("B8").Value = "ABCD867566"
$ws2.Range("C8").Value = "'2190038013"

$ws2.Range("A9").Value = "'2632658757"
$ws2.Range("B9").Value = "ABCD560266"
$ws2.Range("C9").Value = "'5787919800"

$ws2.Range("A10").Value = "'5317762166"
$ws2.Range("B10").Value = "ABCD903584"
$ws2.Range("C10").Value = "'9994515184"

# --- Pick up the freshly prepared FuncLoc / SAID / Previous Doc triple
#     for the live test row (Sheet1 row 2: AV=FuncLoc, AW=SAID,
#     AX=Previous Doc) ------------------------------------------------
$ws.Range("AV2").Value = "ABCD903584"
$ws.Range("AW2").Value = "'9994515184"
$ws.Range("AX2").Value = "5317762166"
